# Update cryptocurrency price/volume data per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.344.35'
$ws.Range('E2').Value = '  -10.54%  '
$ws.Range('D3').Value = '2.306.27'
$ws.Range('E3').Value = '  -20.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '453.02'
$ws.Range('E5').Value = '  -14.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '130.07'
$ws.Range('E6').Value = '  -10.41%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.477'
$ws.Range('E8').Value = '  -14.17%  '
$ws.Range('D9').Value = '2.286.88'
$ws.Range('E9').Value = '  -21.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '5.40'
$ws.Range('E10').Value = '  -10.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').Value = '  -14.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '0.313'
$ws.Range('E12').Value = '  -14.54%  '
$ws.Range('E13').Value = '  -3.06%  '
$ws.Range('D14').Value = '2.677.69'
$ws.Range('E14').Value = '  -21.66%  '
$ws.Range('D15').Value = '54.340.96'
$ws.Range('E15').Value = '  -10.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '18.96'
$ws.Range('E16').Value = '  -16.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.0000121'
$ws.Range('E17').Value = '  -14.83%  '
$ws.Range('D18').Value = '2.268.36'
$ws.Range('E18').Value = '  -22.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '4.13'
$ws.Range('E19').Value = '  -18.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '305.12'
$ws.Range('E20').Value = '  -15.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '9.56'
$ws.Range('E21').Value = '  -18.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '5.60'
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '5.37'
$ws.Range('E24').Value = '  -19.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '56.03'
$ws.Range('E25').Value = '  -13.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '0.989'
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '0.160'
$ws.Range('E27').Value = '  -12.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '0.377'
$ws.Range('E28').Value = '  -17.28%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '6.92'
$ws.Range('E29').Value = '  -12.01%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '0.994'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').Value = '0.0₃0715'
$ws.Range('E31').Value = '  -17.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '145.04'
$ws.Range('E32').Value = '  -3.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '17.09'
$ws.Range('E33').Value = '  -13.65%  '
$ws.Range('E34').Value = '  -19.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '4.78'
$ws.Range('E35').Value = '  -14.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '3.67'
$ws.Range('E36').Value = '  -16.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.849'
$ws.Range('E37').Value = '  -16.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '1.02'
$ws.Range('E38').Value = '  -15.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.991'
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '33.02'
$ws.Range('E40').Value = '  -12.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '10.31'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '1.26'
$ws.Range('E42').Value = '  -15.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '3.19'
$ws.Range('E43').Value = '  -14.77%  '
$ws.Range('D44').Value = '1.942.42'
$ws.Range('E44').Value = '  -15.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.0503'
$ws.Range('E45').Value = '  -13.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.0208'
$ws.Range('E46').Value = '  -12.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.515'
$ws.Range('E47').Value = '  -20.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '0.0819'
$ws.Range('E48').Value = '  -11.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '16.51'
$ws.Range('E49').Value = '  -20.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '4.14'
$ws.Range('E50').Value = '  -18.30%  '
$ws.Range('E51').Value = '  -3.20%  '
